# Applies crypto price/volume updates per commit "Updated cryptos list on Fri Jan  5 09:22:04 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '44.199.01'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.71%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.261.13'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.84%  '

$ws.Range('E4').Value = '  +0.14%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '319.39'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.41%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.42'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.79%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.579'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('E8').Value = '  +0.13%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.555'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.46%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.51'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.01%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0839'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.56%  '

$ws.Range('E12').Value = '  +0.26%  '

$ws.Range('E13').Value = '  -0.77%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.599.56'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.68%  '

$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.861'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.23%  '

$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.48'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.93%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.261.87'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.85%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '44.101.95'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.74%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.33'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.79%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0989'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.58%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.52'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.05%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.80'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.15%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.15'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.58%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '236.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.17%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.11'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.99%  '

$ws.Range('E26').Value = '  +0.07%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.27'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.92%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.19'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.00%  '

$ws.Range('E29').Value = '  -1.95%  '

$ws.Range('E30').Value = '  -1.54%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '161.99'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.56%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.25'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.25%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0853'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.20%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.69'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.47%  '

$ws.Range('E35').Value = '  +11.16%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.98'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.72%  '

$ws.Range('E37').Value = '  -5.88%  '

$ws.Range('B38').Value = 'Stellar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.120'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.60%  '

$ws.Range('B39').Value = 'Celestia'
$ws.Range('C39').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '16.79'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +21.41%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.73'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.43%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.21'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.22%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0317'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.28%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.01'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.39%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.790.31'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.00%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '76.29'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.97%  '

$ws.Range('E46').Value = '  -1.07%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '82.75'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.21%  '

$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.23'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.53%  '

$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '105.32'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.54%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.70'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.82%  '

$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '58.56'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.16%  '
